# Changes of New Pre-Prod URL
# Updates the ShipmentTracking numbers in column P (rows 2-26) of Sheet1
# to the new values produced against the new Pre-Prod URL, while keeping
# the cells stored as text (matching the original shared-string / text type).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$trackingNumbers = @{
    2  = "320018589548"
    3  = "320018589559"
    4  = "320018589581"
    5  = "320018589607"
    6  = "320018589640"
    7  = "320018589662"
    8  = "320018589695"
    9  = "320018589710"
    10 = "320018589743"
    11 = "320018589765"
    12 = "320018589802"
    13 = "320018589824"
    14 = "320018589857"
    15 = "320018589879"
    16 = "320018589905"
    17 = "320018589927"
    18 = "320018589960"
    19 = "320018589982"
    20 = "320018590015"
    21 = "320018590037"
    22 = "320018590060"
    23 = "320018590070"
    24 = "320018590081"
    25 = "320018590092"
    26 = "320018590107"
}

foreach ($row in $trackingNumbers.Keys) {
    $cell = $ws.Range("P$row")
    # Prefix with an apostrophe so the purely-numeric tracking number is
    # stored as text rather than being auto-converted to a number.
    $cell.Value = "'" + $trackingNumbers[$row]
    # Restore the default "Normal" style so no visible formatting change
    # (like the quote-prefix text style) is left on the cell itself.
    $cell.Style = "Normal"
}
